# "add option to leave demand stable"
#
# Inserts a new parameter row on the "Coupling Parameters" sheet (just
# above the existing "yearly_CO2_prices" row, i.e. the new row 17) for a
# "fix_demand_to_initial_year" boolean switch, pushing all subsequent rows
# (and the danger-check formulas that reference them) down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Coupling Parameters")
$ws.Activate()

# Insert a new row above the current row 17 ("yearly_CO2_prices"); Excel
# shifts the rows below down and auto-adjusts the B19/B18-referencing
# formulas to B20/B19 etc.
$ws.Rows("17:17").Insert()

# Fill in the new parameter row. Set the description (C) before the name
# (A) so new shared strings are registered in the same order as the
# authored workbook.
$ws.Range("C17").Value = "so far this is only for NL. DE don't have more than one demand"
$ws.Range("A17").Value = "fix_demand_to_initial_year"
$ws.Range("B17").Value = $true

# Restore/update the view state (scroll position + active cell) like the
# author left it.
$win = $excel.ActiveWindow
$win.ScrollRow = 7
$win.ScrollColumn = 1
$ws.Range("B12").Select()
